$d = $word.ActiveDocument

# Locate the paragraph that starts with "Przygotowano następujące zbiory:"
$target = "Przygotowano następujące zbiory:"
$rng = $d.Content
$found = $rng.Find.Execute($target, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

if (-not $found) {
    throw "Could not find target text: $target"
}

$matchStart = $rng.Start
$matchLen   = $rng.End - $rng.Start

# Build the replacement runs (with proofErr spell-check markers around
# "random.uniform", split exactly as in the target OOXML) as an OOXML
# WordprocessingML package fragment suitable for Range.InsertXML.
$run1 = "Przy pomocy funkcji "
$run2 = "random.uni"
$run3 = "form"
$run4 = "() p"
$run5 = "rzygotowano następujące zbiory:"

$xmlFragment = @"
<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:t xml:space="preserve">$run1</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>$run2</w:t></w:r><w:r><w:t>$run3</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>$run4</w:t></w:r><w:r><w:t>$run5</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
"@

# Insert the new runs right before the matched text (this engine's
# InsertXML merges a bare, attribute-less <w:p> fragment into the
# paragraph addressed by the collapsed range, keeping that paragraph's
# <w:pPr> intact).
$insPoint = $d.Range($matchStart, $matchStart)
$null = $insPoint.InsertXML($xmlFragment)

# Work out how many characters were just inserted, then remove the
# now-shifted original text that used to read "Przygotowano następujące
# zbiory:" immediately after it.
$insertedLen = ($run1 + $run2 + $run3 + $run4 + $run5).Length
$oldRange = $d.Range($matchStart + $insertedLen, $matchStart + $insertedLen + $matchLen)
$oldRange.Delete()
